# Update the example price values on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 0.5
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 1

# Move the active selection to B6 (as left after entering the example data)
[void]$ws.Range("B6").Select()
